# Generate Report for Handback
# - Marks the "6a480caf..." row as failed handback transform (Status column)
#   on the Overview sheet as well as the per-locale (zh-cn / de-de) sheets,
#   since they all share the same underlying text.
# - Fills in the "Error Detail" column for both locale sheets with the
#   specific handback/handoff file-name mismatch message.
# - Widens the "Error Detail" column so the message is readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Update the Status text for the 6a480caf... row everywhere it appears.
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Populate the Error Detail column (P) with the handback/handoff mismatch
# message for each locale.
$zhcn.Range("P3").Value = "Handback file name: 3xhaq2vw.20z is different with handoff file name: 6a480caf-1beb-40a8-aa4a-77bf1c6154ff.680e6086d0c70e827dcb0496a05b83b774f42ccc.zh-cn."
$dede.Range("P3").Value = "Handback file name: 3xhaq2vw.20z is different with handoff file name: 6a480caf-1beb-40a8-aa4a-77bf1c6154ff.680e6086d0c70e827dcb0496a05b83b774f42ccc.de-de."

# Widen the Error Detail column (column 16 / P) on both locale sheets so the
# new longer messages are visible. ColumnWidth uses "characters" units which
# get stored in the XML with a fixed +0.8333 padding offset, so request
# 39.1666... to land exactly on a stored width of 40.
$targetColumnWidth = 39.1666666666667
$zhcn.Columns.Item(16).ColumnWidth = $targetColumnWidth
$dede.Columns.Item(16).ColumnWidth = $targetColumnWidth
